$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 6662.1763
$ws.Range("I6").Value = 919.7692
$ws.Range("J6").Value = 25325
$ws.Range("K6").Value = 2759.3076
$ws.Range("L6").Value = 75975
$ws.Range("M6").Value = -2647.3076
$ws.Range("N6").Value = -76199

$ws.Range("H9").Value = 103.52631
$ws.Range("I9").Value = 111.46667
$ws.Range("J9").Value = 73.75
$ws.Range("K9").Value = 111.46667
$ws.Range("L9").Value = 73.75
$ws.Range("M9").Value = 57.53333000000001
$ws.Range("N9").Value = -411.75

$ws.Range("H12").Value = 170
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""

$ws.Range("H21").Value = 54959.5
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 39900
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 39900
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -40836

$ws.Range("H23").Value = 54959.5
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 39900
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 39900
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -40368

$ws.Range("H29").Value = 824.75
$ws.Range("I29").Value = 433
$ws.Range("K29").Value = 1299
$ws.Range("M29").Value = -1018

$ws.Range("H38").Value = 557.6667
$ws.Range("I38").Value = 557.6667
$ws.Range("K38").Value = 1673.0001
$ws.Range("M38").Value = -1301.0001

$ws.Range("H112").Value = 2020
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 2140
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 6420
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -8636

$ws.Range("H129").Value = 1147.9
$ws.Range("I129").Value = 5398.5
$ws.Range("J129").Value = 844.2857
$ws.Range("K129").Value = 16195.5
$ws.Range("L129").Value = 2532.8571
$ws.Range("M129").Value = -11195.5
$ws.Range("N129").Value = -12532.8571

$ws.Range("H137").Value = 1795830.9
$ws.Range("I137").Value = 7005075
$ws.Range("J137").Value = 5153.1875
$ws.Range("K137").Value = 21015225
$ws.Range("L137").Value = 15459.5625
$ws.Range("M137").Value = -21012675
$ws.Range("N137").Value = -20559.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11518.63
$ws.Range("I32").Value = 11800.391
$ws.Range("K32").Value = 11800.391
$ws.Range("M32").Value = -11513.391

$ws.Range("H88").Value = 2021091.5
$ws.Range("I88").Value = 4625
$ws.Range("J88").Value = 2693247
$ws.Range("K88").Value = 4625
$ws.Range("L88").Value = 2693247
$ws.Range("M88").Value = -4219
$ws.Range("N88").Value = -2694059

$ws.Range("H91").Value = 2021091.5
$ws.Range("I91").Value = 4625
$ws.Range("J91").Value = 2693247
$ws.Range("K91").Value = 4625
$ws.Range("L91").Value = 2693247
$ws.Range("M91").Value = -3221
$ws.Range("N91").Value = -2696055

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 757.85
$ws.Range("I64").Value = 641.4
$ws.Range("J64").Value = 874.3
$ws.Range("K64").Value = 641.4
$ws.Range("L64").Value = 874.3
$ws.Range("M64").Value = -416.4
$ws.Range("N64").Value = -1324.3

$ws.Range("H67").Value = 757.85
$ws.Range("I67").Value = 641.4
$ws.Range("J67").Value = 874.3
$ws.Range("K67").Value = 641.4
$ws.Range("L67").Value = 874.3
$ws.Range("M67").Value = 138.6
$ws.Range("N67").Value = -2434.3

$ws.Range("H115").Value = 20097.572
$ws.Range("J115").Value = 20097.572
$ws.Range("L115").Value = 20097.572
$ws.Range("N115").Value = -23231.572

$ws.Range("H135").Value = 68473.08
$ws.Range("J135").Value = 68473.08
$ws.Range("L135").Value = 68473.08
$ws.Range("N135").Value = -78613.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1202.5385
$ws.Range("I16").Value = 1331.4445
$ws.Range("J16").Value = 912.5
$ws.Range("K16").Value = 1331.4445
$ws.Range("L16").Value = 912.5
$ws.Range("M16").Value = -1044.4445
$ws.Range("N16").Value = -1486.5

$ws.Range("H31").Value = 1718.1718
$ws.Range("I31").Value = 699.383
$ws.Range("J31").Value = 2639
$ws.Range("K31").Value = 699.383
$ws.Range("L31").Value = 2639
$ws.Range("M31").Value = -404.383
$ws.Range("N31").Value = -3229

$ws.Range("H34").Value = 1718.1718
$ws.Range("I34").Value = 699.383
$ws.Range("J34").Value = 2639
$ws.Range("K34").Value = 699.383
$ws.Range("L34").Value = 2639
$ws.Range("M34").Value = -497.383
$ws.Range("N34").Value = -3043

$ws.Range("H52").Value = 53000
$ws.Range("J52").Value = 53000
$ws.Range("L52").Value = 53000
$ws.Range("N52").Value = -53588

$ws.Range("H113").Value = 1202.5385
$ws.Range("I113").Value = 1331.4445
$ws.Range("J113").Value = 912.5
$ws.Range("K113").Value = 1331.4445
$ws.Range("L113").Value = 912.5
$ws.Range("M113").Value = 838.5554999999999
$ws.Range("N113").Value = -5252.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2343.9092
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 2548.3
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 7644.900000000001
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -7812.900000000001

$ws.Range("H39").Value = 4200
$ws.Range("J39").Value = 4200
$ws.Range("L39").Value = 12600
$ws.Range("N39").Value = -13188

$ws.Range("H55").Value = 2499.3845
$ws.Range("J55").Value = 2499.3845
$ws.Range("L55").Value = 7498.1535
$ws.Range("N55").Value = -7852.1535

$ws.Range("H68").Value = 1234
$ws.Range("J68").Value = 1286.9032
$ws.Range("L68").Value = 3860.7096
$ws.Range("N68").Value = -5482.7096

$ws.Range("H71").Value = 1234
$ws.Range("J71").Value = 1286.9032
$ws.Range("L71").Value = 11582.1288
$ws.Range("N71").Value = -19694.1288

$ws.Range("H107").Value = 4052.8035
$ws.Range("I107").Value = 2970.2163
$ws.Range("K107").Value = 8910.6489
$ws.Range("M107").Value = -6990.6489

$ws.Range("H131").Value = 3510.3696
$ws.Range("I131").Value = 11577.667
$ws.Range("J131").Value = 1548.0541
$ws.Range("K131").Value = 34733.001
$ws.Range("L131").Value = 4644.1623
$ws.Range("M131").Value = -29693.001
$ws.Range("N131").Value = -14724.1623

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5405.9116
$ws.Range("I70").Value = 5386.6206
$ws.Range("J70").Value = 5517.8
$ws.Range("K70").Value = 5386.6206
$ws.Range("L70").Value = 5517.8
$ws.Range("M70").Value = -5116.6206
$ws.Range("N70").Value = -6057.8

$ws.Range("H73").Value = 5405.9116
$ws.Range("I73").Value = 5386.6206
$ws.Range("J73").Value = 5517.8
$ws.Range("K73").Value = 5386.6206
$ws.Range("L73").Value = 5517.8
$ws.Range("M73").Value = -4450.6206
$ws.Range("N73").Value = -7389.8

$ws.Range("H97").Value = 11166.833
$ws.Range("I97").Value = 7997.5
$ws.Range("J97").Value = 17505.5
$ws.Range("K97").Value = 7997.5
$ws.Range("L97").Value = 17505.5
$ws.Range("M97").Value = -7501.5
$ws.Range("N97").Value = -18497.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3534.0588
$ws.Range("I100").Value = 2919.889
$ws.Range("J100").Value = 4225
$ws.Range("K100").Value = 2919.889
$ws.Range("L100").Value = 4225
$ws.Range("M100").Value = -2378.889
$ws.Range("N100").Value = -5307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4000
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2627
$ws.Range("N96").Value = ""

$ws.Range("H119").Value = 44499.125
$ws.Range("J119").Value = 44499.125
$ws.Range("L119").Value = 44499.125
$ws.Range("N119").Value = -54175.125

$ws.Range("H132").Value = 1146113.2
$ws.Range("I132").Value = 2289858.2
$ws.Range("J132").Value = 2368.2104
$ws.Range("K132").Value = 6869574.600000001
$ws.Range("L132").Value = 7104.6312
$ws.Range("M132").Value = -6867044.600000001
$ws.Range("N132").Value = -12164.6312
